# Adding logs to the extent report with BasePage and pages
# Edits:
#  - TestCases sheet: widen column A, deselect the tab
#  - TestData sheet: insert a new test-data row (row 5), add its hyperlink,
#    make TestData the active/selected sheet with A5 selected

$wb = $excel.ActiveWorkbook
$wsTestCases = $wb.Worksheets.Item(1)
$wsTestData  = $wb.Worksheets.Item(2)

# --- TestCases (sheet1): widen column A ---
# Excel quantises ColumnWidth to whole pixels (MDW=6 for this workbook's
# default font), so the nearest achievable width to 15.65 chars is reached
# by requesting 14.833333333333334 (-> stored width 15.666666666666666).
$wsTestCases.Columns.Item(1).ColumnWidth = 14.833333333333334

# --- TestData (sheet2): insert new row 5 with test data ---
$wsTestData.Rows.Item(5).Insert()
$wsTestData.Rows.Item(5).RowHeight = 15

$wsTestData.Range("A5").Value = "N"
$wsTestData.Range("B5").Value = "trainer@way2automation.com"
$wsTestData.Range("C5").Value = "askjdfhjskfs"
$wsTestData.Range("D5").Value = "chrome"

# Add the mailto hyperlink on B5 (mirrors B3's hyperlink).
$wsTestData.Hyperlinks.Add($wsTestData.Range("B5"), "mailto:trainer@way2automation.com", "", "", "trainer@way2automation.com")

# Hyperlinks.Add stamps a fresh "visited/hyperlink" style; re-apply B3's
# font so B5 collapses back onto the same shared style B3/B4 already use.
$wsTestData.Range("B5").Font.Underline = $wsTestData.Range("B3").Font.Underline
$wsTestData.Range("B5").Font.Color = $wsTestData.Range("B3").Font.Color
$wsTestData.Range("B5").Font.Name = $wsTestData.Range("B3").Font.Name
$wsTestData.Range("B5").Font.Size = $wsTestData.Range("B3").Font.Size

# --- Active sheet / selection bookkeeping ---
$wsTestData.Activate() | Out-Null
$wsTestData.Range("A5").Select() | Out-Null
